# edit.ps1 - applies the changes described by the diff:
#   1. Update the cached "auto-update" date/time field text from
#      10/29/2019 -> 1/16/2020 everywhere it is cached (slide master +
#      every slide layout's Date placeholder).
#   2. Split the title run on slide 17 ("Scenario 12: ... (watt priority) ")
#      into three runs so that "watt" becomes "var":
#         "Scenario 12: Daily with pf=0.9 and kVA limitation "
#         "(var "
#         "priority) "

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Refresh the cached date field text on the slide master and on every
#    slide layout (PowerPoint stamps this text automatically whenever the
#    deck is saved on a later date; here we push the new cached value).
# ---------------------------------------------------------------------
function Update-DateFieldText {
    param($shapes)

    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.HasTextFrame -and $shp.PlaceholderFormat.Type -eq 16) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "10/29/2019") {
                $tr.Text = "1/16/2020"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateFieldText $master.Shapes

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DateFieldText $layout.Shapes
}

# ---------------------------------------------------------------------
# 2. Slide 17 title: "(watt priority)" -> "(var priority)", landing the
#    new text in its own run(s) just like the authored edit.
# ---------------------------------------------------------------------
$slide17 = $p.Slides.Item(17)
$title = $slide17.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$fullText = $titleRange.Text

$wattIdx = $fullText.IndexOf("(watt ")
if ($wattIdx -ge 0) {
    # Replace "(watt " with "(var " -- becomes its own run.
    $wattRange = $titleRange.Characters($wattIdx + 1, 6)
    $wattRange.Text = "(var "

    # Replace the trailing "priority) " with itself so it becomes a
    # separate run after the one above.
    $afterText = $titleRange.Text
    $priorityIdx = $afterText.IndexOf("priority) ")
    $priorityRange = $titleRange.Characters($priorityIdx + 1, 10)
    $priorityRange.Text = "priority) "
}
